# TestNg Results added, UserDefinedFunctions class
#
# Sheet1 originally has 4 columns: Test Case Name | UserName | Password | Results
# A new "Browser" column is inserted before the existing "Results" column, so
# "Results" (and its "Pass" value) shifts from column D to column E, and the
# new "Browser" column occupies D with the same value ("Pass") as its
# neighboring Results cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column D - "Results" shifts from D to E.
$ws.Columns("D:D").Insert()

# Header for the new column, styled like the rest of the header row.
$ws.Range("D1").Value = "Browser"

# Value for the new column's data row.
$ws.Range("D2").Value = "Pass"

# The source workbook explicitly carries a per-cell style on D2 (matching the
# style already used by A2:C2) rather than leaving it implicit via the column
# default. Toggling MergeCells (a no-op here, nothing is merged) is enough to
# make Excel stamp that explicit style record without altering appearance.
$ws.Range("D2").MergeCells = $false

# Restore the selection to where the editor's cursor ended up.
$ws.Range("D4").Select()
